$d = $word.ActiveDocument

# Locate the last paragraph of body text before the (hidden) bookmarks that
# precede the bibliography ('refs' / 'ref-arrhenius_species_1921'). Those two
# bookmarks are pinned immediately after this paragraph and immediately before
# the following one, so new content must be spliced INTO this paragraph (just
# before its paragraph mark) and then split out into its own paragraphs -
# inserting after/before the paragraph boundary directly would land the new
# text on the wrong side of the bookmarks.
$rngAnchor = $d.Content
$null = $rngAnchor.Find.Execute("lepidoptera, odonata, large mammals", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchorPara = $rngAnchor.Paragraphs(1)
$splicePos = $anchorPara.Range.End - 1
$spliceRange = $d.Range($splicePos, $splicePos)
$spliceRange.InsertAfter("6`tCooperationsThis project takes part in the broader research project of Dr. Petr Keil who has been working on the problem of scale-dependent biodiversity change and integration of heterogeneous data for a decade now, and who has published several high-profile publications on these topics. He currently is my PhD supervisor. Petr’s expertise will be particularly relevant for tasks requiring advanced statistical modelling, interpretation of the models, and putting the results in a broader macroecological context.Cooperation is already ongoing with Vladimír Bejček and Karel Šťastný who furnished us time series of avian biodiversity from the Česka Společnost Ornithologiká which were used in the publications of several atlases (see Bejček and Stastný 2016). On the other hand, Dr. Jiří Reif forwarded me local time series from the Jednotný Program Sčítání Ptáků (JPSP). Their expertise on bird ecologie will be helpful in order to interpret and enhance the outputs of my models.Finally, discussions with Dr. Marta Jarzyna (University of Colombus, Ohio) are ongoing in order to work together on applying the methods that I use on Czech Republic and Europe to some American states.As a matter of fact, my results will 1) allow to differentiate biodiversity dynamics on the North American and the European continents and 2) help to better understand the link between spatio-temporal scales and biodiveristy dynamic by enlarging the databases that I use.")

# --- split the blob into 4 distinct paragraphs, in document order ---
$rh = $d.Content
$null = $rh.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sph = $d.Range($rh.Start, $rh.Start)
$sph.InsertParagraphBefore()
$rp1 = $d.Content
$null = $rp1.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spp1 = $d.Range($rp1.Start, $rp1.Start)
$spp1.InsertParagraphBefore()
$rp2 = $d.Content
$null = $rp2.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spp2 = $d.Range($rp2.Start, $rp2.Start)
$spp2.InsertParagraphBefore()
$rp3 = $d.Content
$null = $rp3.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$spp3 = $d.Range($rp3.Start, $rp3.Start)
$spp3.InsertParagraphBefore()

# --- assign paragraph styles ---
$rs = $d.Content
$null = $rs.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ps = $rs.Paragraphs(1)
$ps.Range.set_Style("Heading1")
$rs = $d.Content
$null = $rs.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ps = $rs.Paragraphs(1)
$ps.Range.set_Style("FirstParagraph")
$rs = $d.Content
$null = $rs.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ps = $rs.Paragraphs(1)
$ps.Range.set_Style("BodyText")
$rs = $d.Content
$null = $rs.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ps = $rs.Paragraphs(1)
$ps.Range.set_Style("BodyText")

# --- bookmark 'coop' wraps the heading text itself ---
$rh2 = $d.Content
$null = $rh2.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headStart = $rh2.Start
$rhEnd = $d.Content
$null = $rhEnd.Find.Execute("6`tCooperations", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$headEnd = $rhEnd.End
$headBookmarkRange = $d.Range($headStart, $headEnd)
$d.Bookmarks.Add("coop", $headBookmarkRange) | Out-Null

# --- hyperlink around '2016' -> internal bookmark ref-bejcek_velke_2016 ---
$rlink = $d.Content
$null = $rlink.Find.Execute("2016", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$linkInner = $d.Range($rlink.Start + 1, $rlink.End - 1)
$d.Hyperlinks.Add($linkInner, [System.Reflection.Missing]::Value, "ref-bejcek_velke_2016", [System.Reflection.Missing]::Value, "2016") | Out-Null

# --- bold the '1)' and '2)' markers ---
$rb = $d.Content
$null = $rb.Find.Execute("1)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boldInner = $d.Range($rb.Start + 1, $rb.End - 1)
$boldInner.Bold = 1
$rb = $d.Content
$null = $rb.Find.Execute("2)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boldInner = $d.Range($rb.Start + 1, $rb.End - 1)
$boldInner.Bold = 1

# --- strip every sentinel marker character left in the document ---
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
$rm = $d.Content
$foundM = $rm.Find.Execute("", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundM) { $rm.Text = "" }
